$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.374.62"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "3.158.87"
$ws.Range("E3").Value = "  +2.26%  "

$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'239.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.42%  "

$ws.Range("D6").Value = "'620.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("D7").Value = "'1.10"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.99%  "

$ws.Range("D8").Value = "'0.372"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.85%  "

$ws.Range("D9").Value = "'0.997"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "

$ws.Range("D10").Value = "'0.743"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.99%  "

$ws.Range("D11").Value = "2.388.32"
$ws.Range("E11").Value = "  -22.62%  "

$ws.Range("E12").Value = "  +2.65%  "

$ws.Range("D13").Value = "'0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.91%  "

$ws.Range("D14").Value = "'35.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("D15").Value = "'5.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.91%  "

$ws.Range("D16").Value = "91.187.68"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("D17").Value = "3.729.70"
$ws.Range("E17").Value = "  +1.63%  "

$ws.Range("D18").Value = "3.151.08"
$ws.Range("E18").Value = "  +1.85%  "

$ws.Range("D19").Value = "'3.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.71%  "

$ws.Range("D20").Value = "'15.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.04%  "

$ws.Range("D21").Value = "'5.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.74%  "

$ws.Range("D22").Value = "'0.0000203"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.01%  "

$ws.Range("D23").Value = "'441.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.72%  "

$ws.Range("D24").Value = "'9.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.21%  "

$ws.Range("D25").Value = "'5.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.41%  "

$ws.Range("D26").Value = "'88.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.85%  "

$ws.Range("D27").Value = "'11.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.81%  "

$ws.Range("D28").Value = "3.314.81"
$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").Value = "'0.124"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +41.06%  "

$ws.Range("D31").Value = "'0.170"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.32%  "

$ws.Range("D32").Value = "'0.228"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +18.63%  "

$ws.Range("D33").Value = "'9.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.88%  "

$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.169"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.05%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "'0.942"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.85%  "

$ws.Range("D36").Value = "'7.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.42%  "

$ws.Range("D37").Value = "'26.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.07%  "

$ws.Range("D38").Value = "'505.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.09%  "

$ws.Range("E39").Value = "  +7.20%  "

$ws.Range("D40").Value = "'1.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.44%  "

$ws.Range("D41").Value = "'0.445"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.40%  "

$ws.Range("D42").Value = "'3.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.02%  "

$ws.Range("D43").Value = "'3.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.85%  "

$ws.Range("D44").Value = "'22.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'0.713"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.51%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'158.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.35%  "

$ws.Range("D48").Value = "'1.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.64%  "

$ws.Range("D49").Value = "'1.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.00%  "

$ws.Range("D50").Value = "'44.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "

$ws.Range("D51").Value = "'4.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.10%  "
